$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row at position 83 - this shifts the existing data rows
# 83..179 down to 84..180 (matching the diff's "dimension A1:T179 -> A1:T180"
# and the observed one-row-down shift of every data row from 83 onward).
$ws.Rows("83:83").Insert()

# Populate the newly inserted row 83 with the new weekly price entry
# (Agrícola del Norte S.A. de Arica - Plátano, fecha 2022-01-10 / serial 44571).
$ws.Cells.Item(83, 1).Value = 1
$ws.Cells.Item(83, 2).Value = "Agrícola del Norte S.A. de Arica"
$ws.Cells.Item(83, 3).Value = "Arica y Parinacota"
$ws.Cells.Item(83, 4).Value = 44571
$ws.Cells.Item(83, 5).Value = 15
$ws.Cells.Item(83, 6).Value = "Fruta"
$ws.Cells.Item(83, 7).Value = 100108
$ws.Cells.Item(83, 8).Value = "Tropicales y subtropicales"
$ws.Cells.Item(83, 9).Value = 100108006
$ws.Cells.Item(83, 10).Value = "Plátano"
$ws.Cells.Item(83, 11).Value = "Sin especificar"
$ws.Cells.Item(83, 12).Value = "Pintón"
$ws.Cells.Item(83, 13).Value = 120
$ws.Cells.Item(83, 14).Value = 14000
$ws.Cells.Item(83, 15).Value = 15000
$ws.Cells.Item(83, 16).Value = 14500
$ws.Cells.Item(83, 17).Value = "$/caja 20 kilos"
$ws.Cells.Item(83, 18).Value = "Bolivia"
$ws.Cells.Item(83, 19).Value = 725
$ws.Cells.Item(83, 20).Value = 20
